# The document starts with two centered paragraphs:
#   1) "AZETECH SOLUTIONS"   (wraps the start of the hidden "_GoBack" bookmark)
#   2) "UART – Level 1 Test"
# and a "_GoBack" bookmark whose end marker sits at the very end of the
# document (after all the body content).
#
# The edit removes the "AZETECH SOLUTIONS" paragraph entirely (text + its
# paragraph mark, merging it into the "UART – Level 1 Test" paragraph) and
# collapses the "_GoBack" bookmark down to an empty bookmark right at the
# start of the document (bookmarkStart immediately followed by
# bookmarkEnd), removing the old trailing bookmarkEnd.

$d = $word.ActiveDocument

# --- Step 1: delete "AZETECH SOLUTIONS" together with its paragraph mark.
# Searching for the literal text plus a carriage-return (paragraph mark)
# character and replacing with nothing merges the two paragraphs, and the
# surviving paragraph keeps its own ("UART – Level 1 Test") formatting.
$d.Content.Find.Execute("AZETECH SOLUTIONS" + [char]13, $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# --- Step 2: move the "_GoBack" bookmark so it collapses right at the
# start of the document instead of wrapping everything through the end of
# the document. Drop the old bookmark first.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Inserting a bookmark directly at position (0,0) is unreliable, so use a
# throw-away placeholder character at the very start, anchor the new
# collapsed bookmark right after it, then delete the placeholder -- the
# bookmark (anchored structurally, not numerically) ends up sitting at
# position 0 once the placeholder is gone.
$placeholder = $d.Range(0, 0)
$placeholder.InsertBefore("X")

$anchor = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $anchor)

$d.Range(0, 1).Delete()
